# "Start and End Time"
# The ManageListings sheet's sample row stored Starttime/Endtime (K2/L2) as
# numeric time-of-day fractions (0.75 / 0.833333...) formatted "hh:mm:ss;@".
# Switch them to plain text time values ("0956pm" / "1144pm") — same text
# style ("@") already used by the neighbouring Startdate/Enddate (H2/I2)
# text cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ManageListings")
$ws.Activate()

$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "0956pm"

$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "1144pm"

# Leave the cursor where the author's session ended up.
$ws.Range("L4").Select()
